$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to Text format so numeric-looking strings
# (e.g. "1.00", "591.44") are preserved exactly as text, matching the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '67.700.47'
$ws.Range('E2').Value = '  +1.16%  '
$ws.Range('D3').Value = '2.533.93'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '591.44'
$ws.Range('E5').Value = '  +0.38%  '
$ws.Range('D6').Value = '172.59'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.525'
$ws.Range('E8').Value = '  -0.54%  '
$ws.Range('D9').Value = '2.533.62'
$ws.Range('E9').Value = '  +0.40%  '
$ws.Range('D10').Value = '0.138'
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('E11').Value = '  +1.27%  '
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('D14').Value = '26.52'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('D15').Value = '0.0000177'
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').Value = '2.919.78'
$ws.Range('E16').Value = '  -2.12%  '
$ws.Range('D17').Value = '67.647.96'
$ws.Range('E17').Value = '  +1.15%  '
$ws.Range('D18').Value = '2.545.98'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').Value = '11.81'
$ws.Range('E19').Value = '  +4.32%  '
$ws.Range('D20').Value = '7.84'
$ws.Range('E20').Value = '  -2.25%  '
$ws.Range('D21').Value = '369.59'
$ws.Range('E21').Value = '  +4.33%  '
$ws.Range('D22').Value = '4.15'
$ws.Range('E22').Value = '  -0.47%  '
$ws.Range('D23').Value = '4.58'
$ws.Range('E23').Value = '  -0.49%  '
$ws.Range('D24').Value = '71.77'
$ws.Range('E24').Value = '  +2.88%  '
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('E26').Value = '  -2.98%  '
$ws.Range('D27').Value = '9.96'
$ws.Range('E27').Value = '  -0.44%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.656.54'
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0964'
$ws.Range('E29').Value = '  -0.98%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '8.43'
$ws.Range('E30').Value = '  +4.27%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').Value = '536.72'
$ws.Range('E31').Value = '  +0.87%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').Value = '1.32'
$ws.Range('E32').Value = '  -0.76%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '1.88'
$ws.Range('E33').Value = '  +2.35%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = '0.129'
$ws.Range('E34').Value = '  -1.36%  '
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').Value = '159.85'
$ws.Range('E36').Value = '  +2.06%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '1.43'
$ws.Range('E37').Value = '  -1.68%  '
$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').Value = '19.09'
$ws.Range('E38').Value = '  +2.80%  '
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').Value = '18.62'
$ws.Range('E39').Value = '  +0.99%  '
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D40').Value = '0.351'
$ws.Range('E40').Value = '  -0.82%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D41').Value = '5.14'
$ws.Range('E41').Value = '  +0.58%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '1.78'
$ws.Range('E42').Value = '  -0.52%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').Value = '2.56'
$ws.Range('E43').Value = '  +3.20%  '
$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '0.0₆0287'
$ws.Range('E45').Value = '  +3.75%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '146.50'
$ws.Range('E46').Value = '  -1.65%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = '3.71'
$ws.Range('E47').Value = '  +0.92%  '
$ws.Range('D48').Value = '0.553'
$ws.Range('E48').Value = '  -0.29%  '
$ws.Range('B49').Value = 'Optimism'
$ws.Range('C49').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D49').Value = '1.72'
$ws.Range('E49').Value = '  +2.15%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.0747'
$ws.Range('E50').Value = '  -1.14%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '0.597'
$ws.Range('E51').Value = '  +0.25%  '
